$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (E1:S1) holds the "dd/mm" headers for each day of attendance.
# The sheet was missing "06/06" entirely and instead had a duplicate/
# mislabeled "16/06" tacked on at the end (column S). Fix this by
# renaming that trailing header to "06/06" and moving it into its
# correct chronological slot (column J, right after "05/06"), which
# shifts the existing "07/06".."15/06" headers one column to the right.
$ws.Range("J1").Value = "06/06"
$ws.Range("K1").Value = "07/06"
$ws.Range("L1").Value = "08/06"
$ws.Range("M1").Value = "09/06"
$ws.Range("N1").Value = "10/06"
$ws.Range("O1").Value = "11/06"
$ws.Range("P1").Value = "12/06"
$ws.Range("Q1").Value = "13/06"
$ws.Range("R1").Value = "14/06"
$ws.Range("S1").Value = "15/06"

# Update the sheet's saved view/selection to match the edited column:
# scroll so column H is the left-most visible column, and select J1.
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
$ws.Range("J1").Select()
